# Append a new data row (31) to Sheet1: Date "2024-10-03", USDValue 2363.64
# The Date column stores plain text (shared string), matching the existing
# rows above it -- so we force text formatting before assigning the value
# to stop Excel's COM layer from auto-converting the "YYYY-MM-DD" string
# into a date serial number, then restore the cell to the sheet's normal
# (unstyled) look so it matches the style-less cells used by every other
# row in that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCell = $ws.Cells.Item(31, 1)
$valueCell = $ws.Cells.Item(31, 2)

$dateCell.NumberFormat = "@"
$dateCell.Value = "2024-10-03"
$dateCell.Style = "Normal"

$valueCell.Value = 2363.64
